# "partial change to step names"
#
# The step-name strings in column D of the "index" sheet are being split
# so that the "T#.#"/"T#" task-reference token no longer carries the
# leading step-number prefix glued onto it: the leading "NN_" and the
# step index that used to sit right after the task token are carved out
# into their own runs (so they can be restyled/re-ordered independently),
# while the surrounding text keeps the same Calibri/11pt/black look.
#
# Net per-cell text result:
#   D2  : 01_1_T2.1_create_conceptset_datasets              -> 01_T2.1_1_create_conceptset_datasets
#   D3  : 01_2_T2.1_create_spells                           -> 01_T2.1_2_create_spells
#   D4  : 01_3_T2.1_create_dates_in_PERSONS                 -> 01_T2.1_3_create_dates_in_PERSONS
#   D5  : (same shared text as D4)                          -> 01_T2.1_3_create_dates_in_PERSONS
#   D6  : 01_4_T2.1_create_prompt_and_itemset_datasets      -> 01_T2.1_4_create_prompt_and_itemset_datasets
#   D7  : (same shared text as D6)                          -> 01_T2.1_4_create_prompt_and_itemset_datasets
#   D8  : 02_1_T2_create_QC_criteria                        -> 02_T2_1_create_QC_criteria
#   D9  : 02_2_T3_apply_QC_exclusion_criteria                -> 02_T3_2_apply_QC_exclusion_criteria
#   D10 : (same shared text as D9)                           -> 02_T3_2_apply_QC_exclusion_criteria
#   D11 : 03_1_T2_create_exclusion_criteria                 -> 03_T2_1_create_exclusion_criteria  (plain text only)
#   D12 : 03_2_T2_merge_persons_concept                      -> 03_T2_2_merge_persons_concept
#   D13 : 04_1_T3_apply_exclusion_criteria                   -> 04_T3_1_apply_exclusion_criteria
#   D14 : (same shared text as D13)                          -> 04_T3_1_apply_exclusion_criteria
#   D15 : (same shared text as D13)                          -> 04_T3_1_apply_exclusion_criteria
#
# NOTE: this engine's PowerShell-function argument binding only works
# reliably with *positional* parameters (`-Name value` named-parameter
# binding silently drops the value), so every helper below is called
# positionally.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("index")

function Set-SplitStep($CellRef, $TaskToken, $StepToken, $Tail) {
    $cell = $ws.Range($CellRef)
    $full = $TaskToken + $StepToken + $Tail
    $cell.Value = $full

    $taskLen = $TaskToken.Length
    $stepLen = $StepToken.Length
    $tailLen = $Tail.Length

    $r1 = $cell.Characters(1, $taskLen)
    $r1.Font.Name = "Calibri"
    $r1.Font.Size = 11
    $r1.Font.Color = 0

    $r2 = $cell.Characters($taskLen + 1, $stepLen)
    $r2.Font.Name = "Calibri"
    $r2.Font.Size = 11
    $r2.Font.Color = 0

    $r3 = $cell.Characters($taskLen + $stepLen + 1, $tailLen)
    $r3.Font.Name = "Calibri"
    $r3.Font.Size = 11
    $r3.Font.Color = 0
}

Set-SplitStep "D2" "01_T2.1" "_1" "_create_conceptset_datasets"
Set-SplitStep "D3" "01_T2.1" "_2" "_create_spells"

Set-SplitStep "D4" "01_T2.1" "_3" "_create_dates_in_PERSONS"
Set-SplitStep "D5" "01_T2.1" "_3" "_create_dates_in_PERSONS"

Set-SplitStep "D6" "01_T2.1" "_4" "_create_prompt_and_itemset_datasets"
Set-SplitStep "D7" "01_T2.1" "_4" "_create_prompt_and_itemset_datasets"

Set-SplitStep "D8" "02_T2" "_1" "_create_QC_criteria"

Set-SplitStep "D9"  "02_T3" "_2" "_apply_QC_exclusion_criteria"
Set-SplitStep "D10" "02_T3" "_2" "_apply_QC_exclusion_criteria"

# D11 is a plain-text rename only (no rich-text run split in the source edit).
$ws.Range("D11").Value = "03_T2_1_create_exclusion_criteria"

Set-SplitStep "D12" "03_T2" "_2" "_merge_persons_concept"

Set-SplitStep "D13" "04_T3" "_1" "_apply_exclusion_criteria"
Set-SplitStep "D14" "04_T3" "_1" "_apply_exclusion_criteria"
Set-SplitStep "D15" "04_T3" "_1" "_apply_exclusion_criteria"
